$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price column (D) to Text format first so numeric-looking values
# (e.g. "0.9996") are not auto-converted to numbers by Excel, keeping
# them as plain text just like the original inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.699.90'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '1.868.01'
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '0.7282'
$ws.Range("E5").Value = '  -0.98%  '

$ws.Range("D6").Value = '240.73'
$ws.Range("E6").Value = '  -0.47%  '

$ws.Range("D7").Value = '0.9992'
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '0.3120'
$ws.Range("E8").Value = '  -1.01%  '

$ws.Range("D9").Value = '0.07081'
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").Value = '0.08224'
$ws.Range("E11").Value = '  -1.88%  '

$ws.Range("D12").Value = '0.7438'
$ws.Range("E12").Value = '  -1.03%  '

$ws.Range("D13").Value = '5.307'
$ws.Range("E13").Value = '  -1.94%  '

$ws.Range("D14").Value = '1.852.05'
$ws.Range("E14").Value = '  -1.00%  '

$ws.Range("D15").Value = '92.10'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '29.696.60'
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").Value = '5.984'
$ws.Range("E17").Value = '  -1.00%  '

$ws.Range("D18").Value = '248.05'
$ws.Range("E18").Value = '  +2.11%  '

$ws.Range("D19").Value = '13.34'
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").Value = '0.000007777'
$ws.Range("E20").Value = '  -0.69%  '

$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").Value = '2.102.72'
$ws.Range("E22").Value = '  -0.63%  '

$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = '7.695'
$ws.Range("E24").Value = '  -2.56%  '

$ws.Range("D25").Value = '0.1529'
$ws.Range("E25").Value = '  -2.43%  '

$ws.Range("D26").Value = '9.148'
$ws.Range("E26").Value = '  -1.74%  '

$ws.Range("D27").Value = '162.84'
$ws.Range("E27").Value = '  -0.81%  '

$ws.Range("D28").Value = '18.48'
$ws.Range("E28").Value = '  -0.49%  '

$ws.Range("D29").Value = '2.005'
$ws.Range("E29").Value = '  -0.52%  '

$ws.Range("E30").Value = '  -2.81%  '

$ws.Range("D31").Value = '4.502'
$ws.Range("E31").Value = '  -2.66%  '

$ws.Range("D32").Value = '1.520'
$ws.Range("E32").Value = '  -0.73%  '

$ws.Range("D33").Value = '4.178'
$ws.Range("E33").Value = '  -2.62%  '

$ws.Range("D34").Value = '0.05263'
$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("D35").Value = '1.227'
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("D36").Value = '0.7513'
$ws.Range("E36").Value = '  -0.22%  '

$ws.Range("D37").Value = '0.9970'
$ws.Range("E37").Value = '  -0.36%  '

$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").Value = '0.01926'
$ws.Range("E39").Value = '  -1.31%  '

$ws.Range("D40").Value = '2.732'
$ws.Range("E40").Value = '  -0.74%  '

$ws.Range("D41").Value = '0.4453'
$ws.Range("E41").Value = '  -0.30%  '

$ws.Range("D42").Value = '5.977'
$ws.Range("E42").Value = '  -1.71%  '

$ws.Range("D43").Value = '0.8651'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").Value = '70.85'
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("D45").Value = '1.044.84'
$ws.Range("E45").Value = '  -5.62%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '103.92'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '0.9998'
$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").Value = '7.451'
$ws.Range("E48").Value = '  -3.27%  '

$ws.Range("D49").Value = '1.812'
$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("D50").Value = '2.004.85'
$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("D51").Value = '2.862'
$ws.Range("E51").Value = '  -6.44%  '

# Restore the Price column style to Normal so no stray style/number-format
# is left referenced on the cells (keeps cells free of an "s" attribute).
$ws.Range("D2:D51").Style = "Normal"
